# Knowledge Enhancement Plan - add a new "Latest status as on 5/23" column (H)
# with status updates for the first two team members, matching the
# formatting already used for the existing "Latest status as on 5/22" (G) column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column H: copy the header/data formatting from column G, then set values.
# (Filled in row order H2, H3, H1 to mirror how the data was authored.)

$ws.Range("G2").Copy()
$ws.Range("H2").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("H2").Value = "1. STAF framework is imported`n2. Installtion process is in progress`n3. Maven errors are being resolved"

$ws.Range("G3").Copy()
$ws.Range("H3").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("H3").Value = "Automated the login page by identifying the page elements into one class and tried to automate the login functionality"

$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("H1").Value = "Latest status as on 5/23"

$excel.CutCopyMode = $false

# --- Resize rows so the new, longer status text is fully visible.
$ws.Rows.Item(1).RowHeight = 45
$ws.Rows.Item(2).RowHeight = 270
$ws.Rows.Item(3).RowHeight = 270

# --- Give the new column a sensible width.
$ws.Columns.Item(8).ColumnWidth = 19.6

# --- Scroll the view back to the left edge of the data and select the new header cell.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 5
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("H1").Select() | Out-Null
